$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the new "English" (C) and "German" (D) translation columns for
#    rows 4-11 first, in row order, so the new shared strings are appended to
#    the shared-strings table in the same order the reference workbook has
#    them (Entertainment, Unterhaltungssektor, Leisure, ... ).
# ---------------------------------------------------------------------------
$translations = @(
    @{ row = 4;  eng = "Entertainment";          deu = "Unterhaltungssektor" },
    @{ row = 5;  eng = "Leisure";                 deu = "Freizeitssektor" },
    @{ row = 6;  eng = "Food";                    deu = "Lebensmittelsektor" },
    @{ row = 7;  eng = "Shopping";                deu = "Einkaufssektor" },
    @{ row = 8;  eng = "Automotive";              deu = "Automobilsektor" },
    @{ row = 9;  eng = "Travel/Tourism";          deu = "Reise-/Tourismussektor" },
    @{ row = 10; eng = "Geographical";            deu = "Geographischer Sektor" },
    @{ row = 11; eng = "Other public services";   deu = "Sonstige öffentliche Dienstleistungen" }
)

foreach ($t in $translations) {
    $ws.Cells.Item($t.row, 3).Value = $t.eng
    $ws.Cells.Item($t.row, 4).Value = $t.deu
}

# Give the new German column (D, rows 4-10) the plain, non-bold "Calibri 12"
# font. D11 is intentionally left with the default style (no explicit font).
$ws.Range("D4:D10").Font.Name = "Calibri"
$ws.Range("D4:D10").Font.Size = 12

# ---------------------------------------------------------------------------
# 2) Replace the old two-column header row (A3:B3 = "Codice"/"Primo Livello")
#    with the new four-column header, and add the header for column C. This
#    is done AFTER the data rows above so "Codice"/"Primo Livello" fall out
#    of use and get pruned, while the four new header strings are appended
#    last in the shared-strings table (matching the target order).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "codice_1_livello"
$ws.Range("B3").Value = "label_ITA_1_livello"
$ws.Range("C3").Value = "label_ENG_1_livello"
$ws.Range("D3").Value = "label_DEU_1_livello"

# Header row keeps the existing bold style already used by A3/B3 (style 1).
$ws.Range("C3").Font.Bold = $true
$ws.Range("D3").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3) Column widths for the new/resized columns.
#    (ColumnWidth is in "characters"; the stored OOXML width attribute is
#    ColumnWidth + 0.8333333333333334 in this engine.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.666666666666668   # -> stored width 20.5
$ws.Columns.Item(3).ColumnWidth = 27.666666666666668   # -> stored width 28.5
$ws.Columns.Item(4).ColumnWidth = 40.998697916666664   # -> stored width ~41.8333 (closest reachable to 41.83203125)

# ---------------------------------------------------------------------------
# 4) Selection moves to C17, matching the saved view state in the target file.
# ---------------------------------------------------------------------------
$ws.Range("C17").Select() | Out-Null
